$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 already holds 2, keep it as-is
$ws.Range("A1").Value = 2

# Header row (order chosen to match the shared-string insertion order)
$ws.Range("K1").Value = "fund"
$ws.Range("O1").Value = "BM Bps Var"
$ws.Range("L1").Value = "NAV per"
$ws.Range("M1").Value = "BM"
$ws.Range("N1").Value = "BM tol"

# Extra "file iteration" columns F/G
$ws.Range("F2").Value = 3123
$ws.Range("G2").Value = 30
$ws.Range("F3").Value = 3478

# Data rows K:O with BM Bps Var formula in O
$ws.Range("K2").Value = 3123
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 130
$ws.Range("N2").Value = 20
$ws.Range("O2").Formula = "=M2-L2"

$ws.Range("K3").Value = 3123
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 130
$ws.Range("N3").Value = 20
$ws.Range("O3").Formula = "=M3-L3"

$ws.Range("K4").Value = 3123
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = 130
$ws.Range("N4").Value = 20
$ws.Range("O4").Formula = "=M4-L4"

$ws.Range("K5").Value = 3456
$ws.Range("L5").Value = 120
$ws.Range("M5").Value = 110
$ws.Range("N5").Value = 40
$ws.Range("O5").Formula = "=M5-L5"

$ws.Range("K6").Value = 3456
$ws.Range("L6").Value = 120
$ws.Range("M6").Value = 110
$ws.Range("N6").Value = 40
$ws.Range("O6").Formula = "=M6-L6"

$ws.Range("K7").Value = 3478
$ws.Range("L7").Value = -60
$ws.Range("M7").Value = -10
$ws.Range("N7").Value = 30
$ws.Range("O7").Formula = "=M7-L7"

$ws.Range("K8").Value = 3478
$ws.Range("L8").Value = -60
$ws.Range("M8").Value = -10
$ws.Range("N8").Value = 30
$ws.Range("O8").Formula = "=M8-L8"

$ws.Range("K9").Value = 3969
$ws.Range("L9").Value = -70
$ws.Range("M9").Value = -60
$ws.Range("N9").Value = 30
$ws.Range("O9").Formula = "=M9-L9"

# Highlight the "BM Bps Var" cells that fall outside the BM tolerance band
$ws.Range("O2").Interior.Color = 65535
$ws.Range("O3").Interior.Color = 65535
$ws.Range("O4").Interior.Color = 65535
$ws.Range("O7").Interior.Color = 65535
$ws.Range("O8").Interior.Color = 65535

# Page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the user last clicked while reviewing results
[void]$ws.Range("H12").Select()
